# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# Adds a new "ODI Bowling Extra" worksheet (mirroring the existing
# "ODI Batting Extra" sheet) populated with MATCH_CODE / MAIDEN_OVERS /
# PERCENT_WICKETS_OF_ALL columns, and tidies up "ODI Batting Extra" by
# dropping the empty placeholder cells that were accidentally written
# for rows that have no batting-extra data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "ODI Bowling Extra" sheet right after
#    "ODI Batting Extra" so tab order matches Player Info / ODI Batting /
#    ODI Bowling / ODI Batting Extra / ODI Bowling Extra.
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$newSheet = $wb.Worksheets.Add($null, $battingExtra)
$newSheet.Name = "ODI Bowling Extra"

# Header row, styled the same way as the other sheets' header rows
# (bold, thin border, centered / top aligned).
$header = $newSheet.Range("A1:C1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "MAIDEN_OVERS"
$newSheet.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Force the data columns to text so numeric-looking strings ("0", "1")
# and percentages ("20.00%") are stored verbatim instead of being
# reinterpreted as numbers / percentages by Excel.
$dataRange = $newSheet.Range("A2:C21")
$dataRange.NumberFormat = "@"

$matchData = @(
    ,@("4453", "0", "20.00%")
    ,@("4455", "0", $null)
    ,@("4608", "1", "10.00%")
    ,@("4614", "0", "30.00%")
    ,@("4625", "0", $null)
    ,@("4636", "0", "10.00%")
    ,@("4639", "0", "10.00%")
    ,@("4642", $null, $null)
    ,@("4647", $null, $null)
    ,@("4648", "1", "10.00%")
    ,@("4649", "0", "10.00%")
    ,@("4669", "0", $null)
    ,@("4673", $null, $null)
    ,@("4676", "0", "10.00%")
    ,@("4686", $null, $null)
    ,@("4688", "0", "10.00%")
    ,@("4690", $null, $null)
    ,@("4692", $null, $null)
    ,@("4695", "0", "10.00%")
    ,@("4697", "0", $null)
)

$rowIndex = 2
foreach ($row in $matchData) {
    $newSheet.Cells.Item($rowIndex, 1).Value = $row[0]

    if ($null -ne $row[1]) {
        $newSheet.Cells.Item($rowIndex, 2).Value = $row[1]
    } else {
        $newSheet.Cells.Item($rowIndex, 2).Value = ""
    }

    if ($null -ne $row[2]) {
        $newSheet.Cells.Item($rowIndex, 3).Value = $row[2]
    } else {
        $newSheet.Cells.Item($rowIndex, 3).Value = ""
    }

    $rowIndex = $rowIndex + 1
}

# ---------------------------------------------------------------------
# 2. Clean up "ODI Batting Extra": remove the empty placeholder cells
#    (NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / BATTING_POSITION) that
#    were written for rows with no data, instead of leaving behind
#    blank-but-present cells.
# ---------------------------------------------------------------------
$emptyCells = @(
    "C2", "D2", "E2",
    "C6", "D6", "E6",
    "B9", "C9", "D9", "E9",
    "B10", "C10", "D10", "E10",
    "C13", "D13", "E13",
    "B14", "C14", "D14", "E14",
    "C15", "D15", "E15",
    "B16", "C16", "D16", "E16",
    "B18", "C18", "D18", "E18",
    "B19", "C19", "D19", "E19"
)

foreach ($cellRef in $emptyCells) {
    $battingExtra.Range($cellRef).Value = $null
}

Write-Output "Added 'ODI Bowling Extra' sheet and tidied 'ODI Batting Extra'."
